# Update the cryptos list sheet with the latest scraped prices / volume
# changes, as produced by the "Updated cryptos list ... with GitHub
# Actions" workflow run.
#
# All cells in the data rows are stored as plain text (the Price column
# routinely holds values Excel would otherwise coerce to a number, e.g.
# "66.698.26" or "7.90"), so every write is forced to text with a leading
# apostrophe and the resulting quote-prefix formatting is stripped with
# ClearFormats() so the cell's style stays untouched (no explicit `s`
# attribute, matching the workbook's original formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.Value = "'" + $value
    $r.ClearFormats()
}

# Rows 42 and 43 swapped rank position (Maker <-> EnergySwap) in addition
# to picking up new price/volume figures, so update B/C/D/E for those two
# rows explicitly.
Set-TextValue "B42" "EnergySwap"
Set-TextValue "C42" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D42" "26.02"
Set-TextValue "E42" "  -1.92%  "

Set-TextValue "B43" "Maker"
Set-TextValue "C43" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D43" "2.792.46"
Set-TextValue "E43" "  +0.13%  "

# Price (D) / Volume(1h) (E) refresh for every other row.
$updates = @{
    2  = @{ D = "66.698.26";  E = "  +0.10%  " }
    3  = @{ D = "3.486.29";   E = "  -0.28%  " }
    4  = @{ E = "  -0.02%  " }
    5  = @{ D = "591.82";     E = "  +0.15%  " }
    6  = @{ D = "171.57";     E = "  +1.53%  " }
    7  = @{ E = "  +0.00%  " }
    8  = @{ D = "0.590";      E = "  -1.71%  " }
    9  = @{ E = "  +2.06%  " }
    10 = @{ D = "7.24";       E = "  -1.36%  " }
    11 = @{ E = "  -2.32%  " }
    12 = @{ D = "4.089.86";   E = "  -0.36%  " }
    13 = @{ E = "  -0.58%  " }
    14 = @{ D = "28.84";      E = "  +2.10%  " }
    15 = @{ D = "66.722.39";  E = "  +0.13%  " }
    16 = @{ E = "  -1.12%  " }
    17 = @{ D = "3.492.95";   E = "  -0.66%  " }
    18 = @{ E = "  -1.20%  " }
    19 = @{ D = "14.05";      E = "  -1.07%  " }
    20 = @{ E = "  -0.34%  " }
    21 = @{ D = "7.90";       E = "  -0.95%  " }
    22 = @{ D = "72.71";      E = "  -0.99%  " }
    23 = @{ E = "  +0.06%  " }
    24 = @{ D = "0.534";      E = "  -0.71%  " }
    25 = @{ E = "  -2.11%  " }
    26 = @{ D = "10.15";      E = "  -0.38%  " }
    27 = @{ E = "  -0.71%  " }
    28 = @{ D = "0.998";      E = "  -0.16%  " }
    29 = @{ E = "  -3.27%  " }
    30 = @{ E = "  -4.16%  " }
    31 = @{ E = "  -1.32%  " }
    32 = @{ D = "23.64";      E = "  +0.04%  " }
    33 = @{ D = "7.31";       E = "  -1.54%  " }
    34 = @{ E = "  -1.28%  " }
    35 = @{ D = "162.93";     E = "  +0.14%  " }
    36 = @{ E = "  -0.73%  " }
    37 = @{ E = "  -1.93%  " }
    38 = @{ D = "6.89";       E = "  +1.88%  " }
    39 = @{ D = "4.65";       E = "  -0.80%  " }
    40 = @{ D = "0.0738";     E = "  -1.15%  " }
    41 = @{ D = "27.16";      E = "  -1.64%  " }
    44 = @{ E = "  -1.29%  " }
    45 = @{ D = "2.54";       E = "  +1.82%  " }
    46 = @{ D = "0.0302";     E = "  -3.52%  " }
    47 = @{ D = "336.61";     E = "  -4.28%  " }
    48 = @{ D = "34.31";      E = "  +1.11%  " }
    49 = @{ E = "  -2.86%  " }
    51 = @{ E = "  -2.46%  " }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        Set-TextValue "$col$row" $cols[$col]
    }
}
